$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.735.46'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.25%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.870.04'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.92%  '

$ws.Range("E4").Value = '  +0.63%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '336.16'
$ws.Range("D5").Style = "Normal"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.011'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.53%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4699'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.39%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3926'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.51%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08014'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.75%  '

$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '45.24'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.11%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.006'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.86'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.67%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.880.33'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.38%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.003'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.26%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.266'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.013'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.41%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.85'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06765'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.08%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001046'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.21%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.29'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.53%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.010'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.44%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.743.71'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.29%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.487'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.68%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.321'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.25%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.098.87'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.52%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.62%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.152'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.78%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.470'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.27%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.97'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.09%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9822'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.35%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09551'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.61%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.633'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.343'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.38%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.343'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.78%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.06068'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.42%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02242'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.33%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.200'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.13%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '8.300'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.42%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.010'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5993'
$ws.Range("D42").Style = "Normal"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1893'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.33'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.30%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.251'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.67%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5677'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.99%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.01%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.933'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.94%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06764'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.10%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '112.47'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.13%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.017'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -11.25%  '
